# This script rewrites the body of the document to match the target revision.
# The "Actual Architecture" milestone doc had a "Database:" paragraph removed,
# several new paragraphs describing PC/hard-drive specifications inserted,
# the "Administration:" paragraph relocated (and its trailing bookmark moved
# up to the title paragraph), and the "Procurement Process:" paragraph text
# reworded. Rather than performing this many fine-grained insert/move/delete
# operations one at a time (error-prone given how much content shifts),
# we replace the whole body content in a single, precise InsertXML call
# that reproduces the exact target OOXML for every paragraph/run.

$d = $word.ActiveDocument

$newBodyXml = @'
<w:p w:rsidR="009971EC" w:rsidRPr="002A33D4" w:rsidRDefault="002A33D4"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="002A33D4"><w:rPr><w:b/></w:rPr><w:t>Actual Architecture</w:t></w:r></w:p><w:p w:rsidR="002A33D4" w:rsidRDefault="002A33D4"><w:r w:rsidRPr="002A33D4"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Client Architecture: </w:t></w:r><w:r><w:t xml:space="preserve"> This system will run on a single Windows based PC.  This PC will be responsible for storing the Access database.</w:t></w:r><w:r><w:t xml:space="preserve">  The PC specifications are:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Brand:</w:t></w:r><w:r><w:t xml:space="preserve"> Dell</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Model:</w:t></w:r><w:r><w:t xml:space="preserve"> I3455-10041WHT</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Memory:</w:t></w:r><w:r><w:t xml:space="preserve"> 8 GB</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Hard Drive:</w:t></w:r><w:r><w:t xml:space="preserve"> 1 TB</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>This specific model was selected because it is an all-in-one unit which makes for a convenient package.  This model also offers sufficient memory for the needs of the system.</w:t></w:r></w:p><w:p><w:r><w:t>The system will also use an external hard drive for back up.  The specifications for this hard drive are:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Brand: </w:t></w:r><w:r><w:t>Western Digital</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/><w:t xml:space="preserve">Model: </w:t></w:r><w:r><w:t>WDBYNN0010BBK-WESN</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:tab/><w:t xml:space="preserve">Hard Drive: </w:t></w:r><w:r><w:t>1 TB</w:t></w:r></w:p><w:p/><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Administration:</w:t></w:r><w:r><w:t xml:space="preserve">  After implementation management staff at The MAX will be responsible for maintenance of this system.  </w:t></w:r></w:p><w:p/><w:p w:rsidR="002A33D4" w:rsidRPr="002A33D4" w:rsidRDefault="002A33D4"><w:r><w:rPr><w:b/></w:rPr><w:t>Procurement Process:</w:t></w:r><w:r><w:t xml:space="preserve">  The MAX </w:t></w:r><w:r><w:t>will purchase this PC</w:t></w:r><w:r><w:t xml:space="preserve"> and hard drive</w:t></w:r><w:r><w:t xml:space="preserve"> from Best Buy.</w:t></w:r></w:p><w:p w:rsidR="002A33D4" w:rsidRDefault="002A33D4"/><w:p w:rsidR="002A33D4" w:rsidRDefault="002A33D4"/><w:p w:rsidR="002A33D4" w:rsidRDefault="002A33D4"/><w:p/>
'@

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Content is the range covering the whole body except the document's final,
# always-present paragraph mark; InsertXML replaces that range's contents
# with the paragraphs encoded above (a trailing empty <w:p/> is included in
# $newBodyXml to account for the paragraph mark that InsertXML itself
# preserves), leaving sectPr and the rest of the package untouched.
$d.Content.InsertXML($pkg)

Write-Host "Body content replaced."
